$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain numeric cells
$ws.Range("A7").Value = 112182733
$ws.Range("B7").Value = 89965
$ws.Range("E7").Value = 760
$ws.Range("Q7").Value = 715207.9695921363
$ws.Range("R7").Value = 7303647.390502339
$ws.Range("S7").Value = 5

# Plain text cells (not mistakable for numbers/dates)
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "VU"
$ws.Range("F7").Value = "Doftticka"
$ws.Range("G7").Value = "Haploporus odorus"
$ws.Range("H7").Value = "(Sommerf.) Bondartsev & Singer"
$ws.Range("P7").Value = "Lomträsk, Pi lm"
$ws.Range("T7").Value = "Norrbotten"
$ws.Range("U7").Value = "Arvidsjaur"
$ws.Range("V7").Value = "Pite lappmark"
$ws.Range("W7").Value = "Arvidsjaur"
$ws.Range("Z7").Value = "00:00"
$ws.Range("AB7").Value = "00:00"
$ws.Range("AC7").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AW7").Value = "Mimmi Persson"
$ws.Range("AX7").Value = "Mimmi Persson"

# Text cells that look like numbers/dates: force text format, then strip
# the formatting override so no stray style survives in styles.xml.
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "1"
$ws.Range("I7").ClearFormats()

$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2023-06-13"
$ws.Range("Y7").ClearFormats()

$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2023-06-13"
$ws.Range("AA7").ClearFormats()

# Booleans
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false

# Blank placeholder cells (present but empty, same as AT6/AY6 in the sheet)
$ws.Range("AT7").NumberFormat = "@"
$ws.Range("AT7").Value = ""
$ws.Range("AT7").ClearFormats()

$ws.Range("AY7").NumberFormat = "@"
$ws.Range("AY7").Value = ""
$ws.Range("AY7").ClearFormats()
